$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(14).Delete()
$ws.Rows(11).Delete()

[void]$ws.Range("A11:B17").Select()
